# Weekly update: a new week of price observations (rows 404-405) is
# inserted at the top of the "Acelga" price history block. Every
# existing weekly pair of rows (Primera/Segunda) from 404-433 shifts
# down by one pair (two rows), so the oldest pair (old 432-433) ends
# up appended as brand-new rows 434-435, and the freshly reported week
# takes over rows 404-405.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing weekly pair (Primera + Segunda rows) down by one
# pair, working from the bottom up so sources aren't overwritten before
# they're copied. This naturally extends the sheet from row 433 to 435
# and carries along every static column (A,B,C,E,F,G,H,I,N,O,Q,R) plus
# the per-row formatting (e.g. the date style on column D).
for ($src = 432; $src -ge 404; $src -= 2) {
    $dest = $src + 2
    $srcRange = $ws.Range("A" + $src + ":R" + ($src + 1))
    $destRange = $ws.Range("A" + $dest)
    $srcRange.Copy($destRange)
}

# Write in the newly reported week's figures on top of rows 404-405
# (everything else on those two rows was already populated by the
# shift above, since row 404/405 used to hold the second-newest week).
$ws.Range("D404").Value = 44746
$ws.Range("J404").Value = 2500
$ws.Range("K404").Value = 600
$ws.Range("L404").Value = 700
$ws.Range("M404").Value = 650
$ws.Range("P404").Value = 325

$ws.Range("D405").Value = 44746
$ws.Range("J405").Value = 1460
$ws.Range("K405").Value = 500
$ws.Range("L405").Value = 550
$ws.Range("M405").Value = 525
$ws.Range("P405").Value = 262
